# Relabel the measure/dimension metadata header table:
#  - row 1: human-readable Spanish labels (was machine slugs)
#  - row 2: iaest-measure / sdmx-dimension identifiers (columns reshuffled;
#           column C becomes "null" and column E becomes the refArea dimension)
#  - row 3: medida/dim role labels for the reshuffled columns
#  - row 4: column C becomes "null"; A/D measures are now xsd:double,
#           B becomes xsd:string, E becomes the URI-Municipio dimension

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - display labels
$ws.Range("A1").Value = "Nº hogares"
$ws.Range("B1").Value = "Número de miembros del hogar"
$ws.Range("C1").Value = "Municipio codigo"
$ws.Range("D1").Value = "Nº medio de miembros"
$ws.Range("E1").Value = "Municipio nombre"

# Row 2 - measure/dimension identifiers
$ws.Range("A2").Value = "iaest-measure:n-hogares"
$ws.Range("B2").Value = "iaest-measure:numero-de-miembros-del-hogar"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "iaest-measure:n-medio-de-miembros"
$ws.Range("E2").Value = "sdmx-dimension:refArea"

# Row 3 - role (medida/dim)
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "dim"

# Row 4 - datatype / URI reference
$ws.Range("A4").Value = "xsd:double"
$ws.Range("B4").Value = "xsd:string"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "xsd:double"
$ws.Range("E4").Value = "URI-Municipio"
